$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> B, C, D, E values (only cells that change are set; unchanged ones left as-is)
$ws.Range("B2").Value  = 6911.9002781658055
$ws.Range("C2").Value  = 1512.9537795510419
$ws.Range("D2").Value  = 4.319937673853629
$ws.Range("E2").Value  = 3.323028979887408

$ws.Range("B3").Value  = 6635.424267039174
$ws.Range("C3").Value  = 1452.4356283690004
$ws.Range("D3").Value  = 4.147140166899485
$ws.Range("E3").Value  = 3.1901078206919116

$ws.Range("B4").Value  = 512.3889750942196
$ws.Range("C4").Value  = 112.15741044730504
$ws.Range("D4").Value  = 0.32024310943388734
$ws.Range("E4").Value  = 0.24634085341068257

$ws.Range("B5").Value  = 14059.713520299201
$ws.Range("C5").Value  = 3077.5468183673474
$ws.Range("D5").Value  = 8.787320950187004
$ws.Range("E5").Value  = 6.759477653990002

$ws.Range("B7").Value  = 6578.612337722757
$ws.Range("D7").Value  = 4.111632711076724
$ws.Range("E7").Value  = 3.162794393135942

$ws.Range("B8").Value  = 1644.6530844306892
$ws.Range("D8").Value  = 1.027908177769181
$ws.Range("E8").Value  = 0.7906985982839855

$ws.Range("B9").Value  = 8223.265422153447
$ws.Range("D9").Value  = 5.1395408888459055
$ws.Range("E9").Value  = 3.9534929914199273

$ws.Range("B11").Value = 3824.4379904672924
$ws.Range("C11").Value = 837.1356182053528
$ws.Range("D11").Value = 2.3902737440420583
$ws.Range("E11").Value = 1.8386721108015833

$ws.Range("C13").Value = 115.39131370412589
$ws.Range("D13").Value = 0.32947687499999995
$ws.Range("E13").Value = 0.25344374999999997

$ws.Range("B14").Value = 1722.5471812266846
$ws.Range("C14").Value = 377.0503281889781
$ws.Range("D14").Value = 1.0765919882666781
$ws.Range("E14").Value = 0.8281476832820601

$ws.Range("C15").Value = 290.24966086707707
$ws.Range("D15").Value = 0.82875
$ws.Range("E15").Value = 0.6375

$ws.Range("C16").Value = 1.1922403470183296
$ws.Range("D16").Value = 0.0034042044515736314
$ws.Range("E16").Value = 0.0026186188089027933

$ws.Range("C17").Value = 10.881200868950028
$ws.Range("D17").Value = 0.031069098214285704
$ws.Range("E17").Value = 0.023899306318681313

$ws.Range("B18").Value = 3438.4529704920587
$ws.Range("C18").Value = 752.6469144741434
$ws.Range("D18").Value = 2.149033106557537
$ws.Range("E18").Value = 1.653102389659644

$ws.Range("B20").Value = 1645.838998047033
$ws.Range("C20").Value = 360.2595859916753
$ws.Range("D20").Value = 1.028649373779396
$ws.Range("E20").Value = 0.7912687490610738

$ws.Range("B21").Value = 5583.811547390809
$ws.Range("C21").Value = 1222.246914008938
$ws.Range("D21").Value = 3.4898822171192565
$ws.Range("E21").Value = 2.684524782399428

$ws.Range("B22").Value = 7350.0095367655185
$ws.Range("C22").Value = 1608.8520176591062
$ws.Range("D22").Value = 4.59375596047845
$ws.Range("E22").Value = 3.5336584311372694

$ws.Range("B25").Value = 36895.879440177516
$ws.Range("C25").Value = 8076.181368705949
$ws.Range("D25").Value = 23.059924650110954
$ws.Range("E25").Value = 17.738403577008423

$ws.Range("B27").Value = 22836.165919878316
$ws.Range("C27").Value = 4998.634550338602
$ws.Range("D27").Value = 14.272603699923952
$ws.Range("E27").Value = 10.978925923018425
